$wb = $excel.ActiveWorkbook

# "Class 2 Formatted" sheet: remove the trailing placeholder row for
# Plumas-Sierra County Fair (FY 04/05) - data not available, shrinking the
# sheet from A1:BV21 to A1:BV20.
$ws1 = $wb.Worksheets.Item("Class 2 Formatted")
$ws1.Rows.Item(21).Delete()

# "Class 5 Formatted" sheet: remove the leading placeholder row for
# 16th DAA, California Mid-State Fair - data not available, shifting the
# remaining fair rows up and shrinking the sheet from A1:BT7 to A1:BT6.
$ws2 = $wb.Worksheets.Item("Class 5 Formatted")
$ws2.Rows.Item(4).Delete()
